$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 47, pushing existing rows 47-67 down to 48-68
$ws.Rows.Item(47).Insert()

# Populate the newly inserted row 47 with the new data record
$ws.Cells.Item(47, 1).Value = 5
$ws.Cells.Item(47, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(47, 3).Value = "Maule"
$ws.Cells.Item(47, 4).Value = 44523
$ws.Cells.Item(47, 4).NumberFormat = $ws.Cells.Item(48, 4).NumberFormat
$ws.Cells.Item(47, 5).Value = 7
$ws.Cells.Item(47, 6).Value = 100112022
$ws.Cells.Item(47, 7).Value = "Arveja Verde"
$ws.Cells.Item(47, 8).Value = "Sin especificar"
$ws.Cells.Item(47, 9).Value = "Primera"
$ws.Cells.Item(47, 10).Value = 600
$ws.Cells.Item(47, 11).Value = 14000
$ws.Cells.Item(47, 12).Value = 14000
$ws.Cells.Item(47, 13).Value = 14000
$ws.Cells.Item(47, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(47, 15).Value = "Región del Maule"
$ws.Cells.Item(47, 16).Value = 560
$ws.Cells.Item(47, 17).Value = 25
$ws.Cells.Item(47, 18).Value = "Hortaliza"
